$wb = $excel.ActiveWorkbook

# Update "F" (想去人数) values on both the "展览" and "全部类型" sheets,
# which mirror the same data.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 517
    $ws.Range("F3").Value = 3451
    $ws.Range("F4").Value = 94
    $ws.Range("F5").Value = 676
}
